$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update title text (July 2021 to Oct 2022 -> July 2021 to Sept 2022)
$ws.Range("A1").Value = "Preventable COVID-19 Deaths, with available vaccination, July 2021 to Sept 2022"

# Update row label (10/2/22 total deaths -> 9/25/22 total deaths)
$ws.Range("A17").Value = "9/25/22 total deaths"

# Update total deaths value for the new date
$ws.Range("B17").Value = 1051501

$wb.Save()
